# Apply the "optimization_parameters" restructuring + active-sheet change
# described by the commit "Updating test files to match the current format
# in beta".

$wb = $excel.ActiveWorkbook

$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsNet = $wb.Worksheets.Item("network_weights")

# --- optimization_parameters sheet content changes -----------------------

# Row 1 (header row): drop the stray duplicate "value" labels in C1:F1 -
# only A1/B1 ("optimization_parameter" / "value") remain.
$wsOpt.Range("C1:F1").ClearContents()

# Row 8, column A: "Model" -> "production_function" (value in B8 unchanged).
$wsOpt.Range("A8").Value = "production_function"

# Insert a brand-new row directly below (becomes the new row 9) to hold the
# "L_curve" parameter, pushing the old rows 9-16 down by one.
$wsOpt.Rows.Item(9).Insert()
$wsOpt.Range("A9").Value = "L_curve"
$wsOpt.Range("B9").Value = 0
$wsOpt.Range("B9").NumberFormat = "0.00E+00"

# The old "Deletion" row (now shifted to row 17) is removed entirely.
$wsOpt.Rows.Item(17).Delete()

# --- active sheet / tab-selection change ----------------------------------
# Previously "network_weights" was the active/selected tab with B2:E5
# selected; now "optimization_parameters" is active, with C1:F1 selected.
$wsOpt.Activate()
$wsOpt.Range("C1:F1").Select()
